$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 67

# Columns A (Date) and D (Week) hold values that look like dates/numbers
# but must stay as literal text, matching the rest of the sheet's rows.
# Force text format before assigning so Excel doesn't auto-convert them.
$ws.Range("A67").NumberFormat = "@"
$ws.Range("D67").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2023-06-22"
$ws.Cells.Item($row, 2).Value = "11:40:18"
$ws.Cells.Item($row, 3).Value = "Thursday"
$ws.Cells.Item($row, 4).Value = "25"

$ws.Cells.Item($row, 5).Value = 122376
$ws.Cells.Item($row, 6).Value = 133742
$ws.Cells.Item($row, 7).Value = 162487
$ws.Cells.Item($row, 8).Value = 133692
$ws.Cells.Item($row, 9).Value = 177333
$ws.Cells.Item($row, 10).Value = 114995
$ws.Cells.Item($row, 11).Value = 202290
$ws.Cells.Item($row, 12).Value = 225743
$ws.Cells.Item($row, 13).Value = 175544
$ws.Cells.Item($row, 14).Value = 104000
$ws.Cells.Item($row, 15).Value = 39392
$ws.Cells.Item($row, 16).Value = 33891
$ws.Cells.Item($row, 17).Value = 51938
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 35952
$ws.Cells.Item($row, 20).Value = -1
